$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F ("Lista 2") attendance/points update after 29.10.2018 ---
# The cells in column F were previously blank placeholders; copy the cell
# formatting (fill/border) used elsewhere in the sheet onto each F cell so the
# visual style matches what Excel would apply when a value is typed in, then
# set the actual values.

# Row 2 uses the "obecność bez wysłanej listy" style (same as E2)
$ws.Range("E2").Copy($ws.Range("F2"))
$ws.Range("F2").Value = 5

# Row 3 uses the regular style (same as E3)
$ws.Range("E3").Copy($ws.Range("F3"))
$ws.Range("F3").Value = 9

# Row 4
$ws.Range("E3").Copy($ws.Range("F4"))
$ws.Range("F4").Value = 6

# Row 5
$ws.Range("E3").Copy($ws.Range("F5"))
$ws.Range("F5").Value = 3

# Row 6
$ws.Range("E3").Copy($ws.Range("F6"))
$ws.Range("F6").Value = 3

# Row 7 uses the "nieobecność usprawiedliwiona" style (same as E7), stays empty
$ws.Range("E7").Copy($ws.Range("F7"))
$ws.Range("F7").ClearContents()

# Row 8
$ws.Range("E3").Copy($ws.Range("F8"))
$ws.Range("F8").Value = 9

# Row 9 stays empty but gets the regular style
$ws.Range("E3").Copy($ws.Range("F9"))
$ws.Range("F9").ClearContents()

# Row 10
$ws.Range("E3").Copy($ws.Range("F10"))
$ws.Range("F10").Value = 6

# Row 11 uses the "nieobecność usprawiedliwiona" style, stays empty
$ws.Range("E7").Copy($ws.Range("F11"))
$ws.Range("F11").ClearContents()

# Row 12
$ws.Range("E3").Copy($ws.Range("F12"))
$ws.Range("F12").Value = 5

# Row 13: attendance count in E13 corrected from 2 to 3, and F13 filled in
$ws.Range("E13").Value = 3
$ws.Range("E3").Copy($ws.Range("F13"))
$ws.Range("F13").Value = 7

# Row 14
$ws.Range("E3").Copy($ws.Range("F14"))
$ws.Range("F14").Value = 3

# Row 15
$ws.Range("E3").Copy($ws.Range("F15"))
$ws.Range("F15").Value = 5

# Row 16 (last row of the table, bottom border) - take style from D16 which
# already carries the regular look for this row
$ws.Range("D16").Copy($ws.Range("F16"))
$ws.Range("F16").Value = 7

# Update the active selection to reflect where the author left off editing
$ws.Range("F15").Select()
